# Update "Estado de Cuenta" worker rows (B15:J21 table): the underlying
# database was refreshed, so the workers' records (doc number, name,
# "Valor Mora" and "Salario Basico") were reshuffled across rows 16-21
# while the template/layout stays the same (Tipo Doc "CC" and Periodo
# Mora "1907" already constant for every row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 16; Doc = "8867002";    Name = "JAMED ANTONIO HOYOS RAMOS";     Mora = 72000;  Salario = 2151340 },
    @{ Row = 17; Doc = "1101813586"; Name = "YESID FABIAN CARDENAS OLIVERA"; Mora = 35616;  Salario = 2000000 },
    @{ Row = 18; Doc = "1002377491"; Name = "WILLIAN JOSE ECHAVEZ MORALES";  Mora = 160000; Salario = 5198000 },
    @{ Row = 19; Doc = "1002372088"; Name = "JORGE ALFREDO PORTELA MARTINEZ"; Mora = 120000; Salario = 3430000 },
    @{ Row = 20; Doc = "1102893944"; Name = "ALEXANDRA ISABEL MENCO MORALES"; Mora = 40000;  Salario = 2260000 },
    @{ Row = 21; Doc = "1045671855"; Name = "JORGE LEONARDO ESTRADA YANES";  Mora = 72000;  Salario = 1991340 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Name
    $ws.Cells.Item($r, 6).Value = $item.Mora
    $ws.Cells.Item($r, 7).Value = $item.Salario
}
